$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) holds values that look numeric (e.g. "580.60", "0.0000169")
# but are stored as text in the workbook. Setting .Value directly would let Excel
# auto-convert these into real numbers (losing trailing zeros / switching to
# scientific notation), so every Price cell we touch is first forced to Text format.
foreach ($addr in @("D2", "D3", "D5", "D6", "D9", "D13", "D15", "D16", "D17", "D18", "D19", "D20", "D21", "D22", "D24", "D27", "D29", "D31", "D32", "D35", "D37", "D38", "D40", "D41", "D44", "D45", "D46", "D47")) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "67.847.95"
$ws.Range("E2").Value = "  +4.43%  "
$ws.Range("D3").Value = "3.271.13"
$ws.Range("E3").Value = "  +4.29%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "580.60"
$ws.Range("E5").Value = "  +2.24%  "
$ws.Range("D6").Value = "182.83"
$ws.Range("E6").Value = "  +8.54%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  +0.17%  "
$ws.Range("D9").Value = "3.270.30"
$ws.Range("E9").Value = "  +4.34%  "
$ws.Range("E10").Value = "  +8.40%  "
$ws.Range("E11").Value = "  +3.69%  "
$ws.Range("E12").Value = "  +7.35%  "
$ws.Range("D13").Value = "3.834.83"
$ws.Range("E13").Value = "  +4.39%  "
$ws.Range("E14").Value = "  +1.20%  "
$ws.Range("D15").Value = "28.66"
$ws.Range("E15").Value = "  +7.53%  "
$ws.Range("D16").Value = "67.787.18"
$ws.Range("E16").Value = "  +4.56%  "
$ws.Range("D17").Value = "0.0000169"
$ws.Range("E17").Value = "  +4.88%  "
$ws.Range("D18").Value = "3.267.60"
$ws.Range("E18").Value = "  +4.53%  "
$ws.Range("D19").Value = "5.86"
$ws.Range("E19").Value = "  +3.24%  "
$ws.Range("D20").Value = "13.63"
$ws.Range("E20").Value = "  +7.35%  "
$ws.Range("D21").Value = "375.85"
$ws.Range("E21").Value = "  +6.07%  "
$ws.Range("D22").Value = "7.67"
$ws.Range("E22").Value = "  +6.35%  "
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("D24").Value = "71.34"
$ws.Range("E24").Value = "  +3.78%  "
$ws.Range("E25").Value = "  +4.35%  "
$ws.Range("E26").Value = "  +5.77%  "
$ws.Range("D27").Value = "9.67"
$ws.Range("E27").Value = "  +0.59%  "
$ws.Range("E28").Value = "  +3.63%  "
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("E30").Value = "  +4.30%  "
$ws.Range("D31").Value = "5.73"
$ws.Range("E31").Value = "  +9.34%  "
$ws.Range("D32").Value = "22.77"
$ws.Range("E32").Value = "  +4.87%  "
$ws.Range("E33").Value = "  -0.13%  "
$ws.Range("E34").Value = "  +8.53%  "
$ws.Range("D35").Value = "6.95"
$ws.Range("E35").Value = "  +6.21%  "
$ws.Range("E36").Value = "  +6.59%  "
$ws.Range("D37").Value = "163.20"
$ws.Range("E37").Value = "  +3.04%  "
$ws.Range("D38").Value = "0.854"
$ws.Range("E38").Value = "  +3.36%  "
$ws.Range("E39").Value = "  +5.97%  "
$ws.Range("D40").Value = "6.86"
$ws.Range("E40").Value = "  +13.04%  "
$ws.Range("D41").Value = "4.69"
$ws.Range("E41").Value = "  +13.16%  "
$ws.Range("E42").Value = "  +3.43%  "
$ws.Range("E43").Value = "  +9.71%  "
$ws.Range("D44").Value = "356.62"
$ws.Range("E44").Value = "  +12.53%  "
$ws.Range("D45").Value = "2.712.31"
$ws.Range("E45").Value = "  +2.75%  "
$ws.Range("D46").Value = "25.52"
$ws.Range("E46").Value = "  +7.00%  "
$ws.Range("D47").Value = "40.94"
$ws.Range("E47").Value = "  +4.15%  "
$ws.Range("E48").Value = "  +5.54%  "
$ws.Range("E49").Value = "  +4.53%  "
$ws.Range("E50").Value = "  +7.94%  "
$ws.Range("E51").Value = "  +1.08%  "
